# Re-pull data / push all data / recompute mean -> updates to the dSF (column F) values
# for a number of rows in the season log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new dSF (column F) value
$updates = @{
    6  = -3
    8  = -1
    9  = -2
    16 = -2
    18 = 3
    22 = -1
    32 = -1
    37 = 4
    41 = 10
    43 = 1
    48 = 2
    53 = 0
    57 = -1
    59 = 0
    69 = 0
    71 = -2
    75 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
